# Reassign the grid_cell (column AG) values for the distr_solelc_won-CHE_00xx
# rows on the "solar" sheet to the updated grid-cell mapping.
$wb = $excel.ActiveWorkbook
$sheet = $wb.Worksheets.Item("solar")

$sheet.Range("AG4").Value = "CHE_14"
$sheet.Range("AG5").Value = "CHE_18"
$sheet.Range("AG6").Value = "CHE_9"
$sheet.Range("AG7").Value = "CHE_21"
$sheet.Range("AG8").Value = "CHE_4"
$sheet.Range("AG9").Value = "CHE_0"
$sheet.Range("AG10").Value = "CHE_11"
$sheet.Range("AG11").Value = "CHE_15"
$sheet.Range("AG12").Value = "CHE_25"
$sheet.Range("AG13").Value = "CHE_2"
$sheet.Range("AG14").Value = "CHE_10"
$sheet.Range("AG15").Value = "CHE_22"
$sheet.Range("AG16").Value = "CHE_17"
$sheet.Range("AG17").Value = "CHE_19"
$sheet.Range("AG18").Value = "CHE_23"
$sheet.Range("AG19").Value = "CHE_7"
$sheet.Range("AG20").Value = "CHE_13"
$sheet.Range("AG21").Value = "CHE_20"
$sheet.Range("AG22").Value = "CHE_1"
$sheet.Range("AG23").Value = "CHE_6"
$sheet.Range("AG24").Value = "CHE_24"
$sheet.Range("AG25").Value = "CHE_8"
$sheet.Range("AG26").Value = "CHE_5"
$sheet.Range("AG27").Value = "CHE_3"
$sheet.Range("AG28").Value = "CHE_12"
